$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 2
$ws.Range("H2").Value = 3.9
$ws.Range("I2").Value = 3.4
$ws.Range("J2").Value = 2.6
$ws.Range("K2").Value = 2.4
$ws.Range("L2").Value = 3.6
$ws.Range("N2").Value = 15
$ws.Range("Q2").Value = 1.57
$ws.Range("R2").Value = 2.38
$ws.Range("S2").Value = 1.29
$ws.Range("T2").Value = 3.5
$ws.Range("AA2").Value = 15
$ws.Range("AB2").Value = 21
$ws.Range("AC2").Value = 17
$ws.Range("AD2").Value = 7.5
$ws.Range("AP2").Value = 17
$ws.Range("AQ2").Value = 34
$ws.Range("AS2").Value = 101
$ws.Range("AT2").Value = 3.5
$ws.Range("AW2").Value = 351
$ws.Range("AX2").Value = 5.5
